{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// The document's first paragraph holds a date string, and a single 20x5\n// table follows it where every cell holds exactly one paragraph / one run\n// with a \"A\u00d7B=\" text. The diff replaces the date plus every cell's text,\n// in document order (paragraph index 0 = date paragraph, 1..100 = the\n// table cells read row-major). One \"old\" value (\"57\u00d715=\") repeats twice\n// with two different replacements, so we must walk the paragraphs in\n// document order (positional) rather than doing a global text search and\n// replace.\n\nconst REPLACEMENTS = [\n  [\"2023-06-30 Friday\", \"2023-07-01 Saturday\"],\n  [\"91\u00d761=\", \"99\u00d747=\"],\n  [\"69\u00d795=\", \"93\u00d774=\"],\n  [\"70\u00d731=\", \"82\u00d792=\"],\n  [\"94\u00d794=\", \"67\u00d777=\"],\n  [\"17\u00d786=\", \"55\u00d767=\"],\n  [\"13\u00d740=\", \"50\u00d765=\"],\n  [\"55\u00d728=\", \"24\u00d7100=\"],\n  [\"28\u00d782=\", \"38\u00d729=\"],\n  [\"32\u00d786=\", \"72\u00d767=\"],\n  [\"32\u00d715=\", \"85\u00d748=\"],\n  [\"47\u00d772=\", \"68\u00d741=\"],\n  [\"95\u00d752=\", \"24\u00d770=\"],\n  [\"31\u00d745=\", \"64\u00d786=\"],\n  [\"81\u00d749=\", \"40\u00d764=\"],\n  [\"29\u00d737=\", \"41\u00d714=\"],\n  [\"88\u00d711=\", \"63\u00d790=\"],\n  [\"95\u00d786=\", \"83\u00d712=\"],\n  [\"100\u00d720=\", \"41\u00d755=\"],\n  [\"82\u00d776=\", \"43\u00d716=\"],\n  [\"12\u00d744=\", \"89\u00d742=\"],\n  [\"73\u00d774=\", \"64\u00d772=\"],\n  [\"57\u00d715=\", \"16\u00d726=\"],\n  [\"38\u00d719=\", \"64\u00d747=\"],\n  [\"76\u00d741=\", \"66\u00d774=\"],\n  [\"96\u00d729=\", \"68\u00d767=\"],\n  [\"52\u00d743=\", \"37\u00d786=\"],\n  [\"84\u00d772=\", \"57\u00d755=\"],\n  [\"12\u00d717=\", \"82\u00d722=\"],\n  [\"80\u00d713=\", \"27\u00d771=\"],\n  [\"77\u00d781=\", \"47\u00d797=\"],\n  [\"76\u00d781=\", \"44\u00d768=\"],\n  [\"39\u00d750=\", \"81\u00d721=\"],\n  [\"44\u00d751=\", \"25\u00d750=\"],\n  [\"70\u00d784=\", \"65\u00d719=\"],\n  [\"71\u00d712=\", \"99\u00d734=\"],\n  [\"17\u00d739=\", \"86\u00d789=\"],\n  [\"90\u00d793=\", \"69\u00d731=\"],\n  [\"46\u00d780=\", \"66\u00d728=\"],\n  [\"97\u00d755=\", \"23\u00d743=\"],\n  [\"57\u00d715=\", \"32\u00d716=\"],\n  [\"90\u00d740=\", \"29\u00d717=\"],\n  [\"95\u00d731=\", \"22\u00d714=\"],\n  [\"10\u00d754=\", \"30\u00d771=\"],\n  [\"38\u00d764=\", \"89\u00d767=\"],\n  [\"36\u00d712=\", \"16\u00d724=\"],\n  [\"14\u00d740=\", \"72\u00d726=\"],\n  [\"74\u00d725=\", \"56\u00d723=\"],\n  [\"69\u00d777=\", \"18\u00d713=\"],\n  [\"98\u00d753=\", \"81\u00d796=\"],\n  [\"24\u00d759=\", \"48\u00d765=\"],\n  [\"42\u00d767=\", \"80\u00d763=\"],\n  [\"49\u00d727=\", \"66\u00d726=\"],\n  [\"52\u00d713=\", \"49\u00d772=\"],\n  [\"69\u00d784=\", \"14\u00d799=\"],\n  [\"46\u00d735=\", \"53\u00d721=\"],\n  [\"49\u00d789=\", \"20\u00d712=\"],\n  [\"91\u00d778=\", \"68\u00d790=\"],\n  [\"35\u00d772=\", \"33\u00d760=\"],\n  [\"85\u00d753=\", \"84\u00d785=\"],\n  [\"72\u00d763=\", \"91\u00d720=\"],\n  [\"74\u00d754=\", \"14\u00d766=\"],\n  [\"46\u00d753=\", \"55\u00d750=\"],\n  [\"76\u00d784=\", \"35\u00d773=\"],\n  [\"90\u00d768=\", \"89\u00d728=\"],\n  [\"85\u00d738=\", \"59\u00d775=\"],\n  [\"46\u00d762=\", \"86\u00d779=\"],\n  [\"68\u00d746=\", \"95\u00d732=\"],\n  [\"96\u00d741=\", \"13\u00d711=\"],\n  [\"32\u00d774=\", \"94\u00d793=\"],\n  [\"21\u00d741=\", \"95\u00d751=\"],\n  [\"90\u00d715=\", \"72\u00d713=\"],\n  [\"20\u00d768=\", \"33\u00d728=\"],\n  [\"100\u00d784=\", \"43\u00d769=\"],\n  [\"17\u00d796=\", \"82\u00d797=\"],\n  [\"58\u00d760=\", \"32\u00d711=\"],\n  [\"59\u00d731=\", \"81\u00d725=\"],\n  [\"10\u00d712=\", \"15\u00d778=\"],\n  [\"24\u00d747=\", \"18\u00d745=\"],\n  [\"24\u00d780=\", \"69\u00d753=\"],\n  [\"28\u00d773=\", \"38\u00d736=\"],\n  [\"36\u00d750=\", \"47\u00d733=\"],\n  [\"54\u00d740=\", \"16\u00d756=\"],\n  [\"83\u00d727=\", \"19\u00d767=\"],\n  [\"61\u00d758=\", \"41\u00d716=\"],\n  [\"97\u00d762=\", \"86\u00d740=\"],\n  [\"46\u00d766=\", \"23\u00d767=\"],\n  [\"23\u00d747=\", \"63\u00d720=\"],\n  [\"64\u00d744=\", \"52\u00d712=\"],\n  [\"22\u00d728=\", \"12\u00d797=\"],\n  [\"68\u00d723=\", \"96\u00d747=\"],\n  [\"46\u00d793=\", \"91\u00d773=\"],\n  [\"98\u00d719=\", \"15\u00d781=\"],\n  [\"10\u00d776=\", \"70\u00d771=\"],\n  [\"12\u00d730=\", \"76\u00d770=\"],\n  [\"53\u00d764=\", \"41\u00d771=\"],\n  [\"35\u00d759=\", \"72\u00d735=\"],\n  [\"16\u00d771=\", \"57\u00d785=\"],\n  [\"18\u00d738=\", \"34\u00d713=\"],\n  [\"60\u00d737=\", \"12\u00d745=\"],\n  [\"81\u00d754=\", \"97\u00d752=\"],\n];\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== REPLACEMENTS.length) {\n  throw new Error(\n    \"Expected \" + REPLACEMENTS.length + \" paragraphs (1 date + 100 table cells), found \" +\n    paragraphs.items.length\n  );\n}\n\n// Load each paragraph's current text so we can sanity-check against the\n// expected \"before\" value prior to overwriting it.\nfor (const p of paragraphs.items) p.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < REPLACEMENTS.length; i++) {\n  const [oldText, newText] = REPLACEMENTS[i];\n  const para = paragraphs.items[i];\n  const current = (para.text || \"\").replace(/\\r$/, \"\");\n  if (current !== oldText) {\n    throw new Error(\n      \"Paragraph \" + i + \": expected \" + JSON.stringify(oldText) +\n      \" but found \" + JSON.stringify(current)\n    );\n  }\n  // Replace just the text, preserving the run's formatting (font, size, ...).\n  para.insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document ($d below).\n#\n# The document starts with one paragraph holding a date string, followed by\n# a single 20-row x 5-column table where every cell holds exactly one\n# paragraph / one run of \"A\u00d7B=\" text. The diff rewrites the date plus every\n# cell's text. One \"old\" cell value (\"57\u00d715=\") occurs twice with two\n# different replacements, so cells are addressed positionally via\n# Tables(1).Cell(row, col) (row-major, matching the document/diff order)\n# rather than via a blanket Find/Replace-All.\n\n$d = $word.ActiveDocument\n\n$OLD_VALUES = @(\n  \"91\u00d761=\",\n  \"69\u00d795=\",\n  \"70\u00d731=\",\n  \"94\u00d794=\",\n  \"17\u00d786=\",\n  \"13\u00d740=\",\n  \"55\u00d728=\",\n  \"28\u00d782=\",\n  \"32\u00d786=\",\n  \"32\u00d715=\",\n  \"47\u00d772=\",\n  \"95\u00d752=\",\n  \"31\u00d745=\",\n  \"81\u00d749=\",\n  \"29\u00d737=\",\n  \"88\u00d711=\",\n  \"95\u00d786=\",\n  \"100\u00d720=\",\n  \"82\u00d776=\",\n  \"12\u00d744=\",\n  \"73\u00d774=\",\n  \"57\u00d715=\",\n  \"38\u00d719=\",\n  \"76\u00d741=\",\n  \"96\u00d729=\",\n  \"52\u00d743=\",\n  \"84\u00d772=\",\n  \"12\u00d717=\",\n  \"80\u00d713=\",\n  \"77\u00d781=\",\n  \"76\u00d781=\",\n  \"39\u00d750=\",\n  \"44\u00d751=\",\n  \"70\u00d784=\",\n  \"71\u00d712=\",\n  \"17\u00d739=\",\n  \"90\u00d793=\",\n  \"46\u00d780=\",\n  \"97\u00d755=\",\n  \"57\u00d715=\",\n  \"90\u00d740=\",\n  \"95\u00d731=\",\n  \"10\u00d754=\",\n  \"38\u00d764=\",\n  \"36\u00d712=\",\n  \"14\u00d740=\",\n  \"74\u00d725=\",\n  \"69\u00d777=\",\n  \"98\u00d753=\",\n  \"24\u00d759=\",\n  \"42\u00d767=\",\n  \"49\u00d727=\",\n  \"52\u00d713=\",\n  \"69\u00d784=\",\n  \"46\u00d735=\",\n  \"49\u00d789=\",\n  \"91\u00d778=\",\n  \"35\u00d772=\",\n  \"85\u00d753=\",\n  \"72\u00d763=\",\n  \"74\u00d754=\",\n  \"46\u00d753=\",\n  \"76\u00d784=\",\n  \"90\u00d768=\",\n  \"85\u00d738=\",\n  \"46\u00d762=\",\n  \"68\u00d746=\",\n  \"96\u00d741=\",\n  \"32\u00d774=\",\n  \"21\u00d741=\",\n  \"90\u00d715=\",\n  \"20\u00d768=\",\n  \"100\u00d784=\",\n  \"17\u00d796=\",\n  \"58\u00d760=\",\n  \"59\u00d731=\",\n  \"10\u00d712=\",\n  \"24\u00d747=\",\n  \"24\u00d780=\",\n  \"28\u00d773=\",\n  \"36\u00d750=\",\n  \"54\u00d740=\",\n  \"83\u00d727=\",\n  \"61\u00d758=\",\n  \"97\u00d762=\",\n  \"46\u00d766=\",\n  \"23\u00d747=\",\n  \"64\u00d744=\",\n  \"22\u00d728=\",\n  \"68\u00d723=\",\n  \"46\u00d793=\",\n  \"98\u00d719=\",\n  \"10\u00d776=\",\n  \"12\u00d730=\",\n  \"53\u00d764=\",\n  \"35\u00d759=\",\n  \"16\u00d771=\",\n  \"18\u00d738=\",\n  \"60\u00d737=\",\n  \"81\u00d754=\",\n)\n\n$NEW_VALUES = @(\n  \"99\u00d747=\",\n  \"93\u00d774=\",\n  \"82\u00d792=\",\n  \"67\u00d777=\",\n  \"55\u00d767=\",\n  \"50\u00d765=\",\n  \"24\u00d7100=\",\n  \"38\u00d729=\",\n  \"72\u00d767=\",\n  \"85\u00d748=\",\n  \"68\u00d741=\",\n  \"24\u00d770=\",\n  \"64\u00d786=\",\n  \"40\u00d764=\",\n  \"41\u00d714=\",\n  \"63\u00d790=\",\n  \"83\u00d712=\",\n  \"41\u00d755=\",\n  \"43\u00d716=\",\n  \"89\u00d742=\",\n  \"64\u00d772=\",\n  \"16\u00d726=\",\n  \"64\u00d747=\",\n  \"66\u00d774=\",\n  \"68\u00d767=\",\n  \"37\u00d786=\",\n  \"57\u00d755=\",\n  \"82\u00d722=\",\n  \"27\u00d771=\",\n  \"47\u00d797=\",\n  \"44\u00d768=\",\n  \"81\u00d721=\",\n  \"25\u00d750=\",\n  \"65\u00d719=\",\n  \"99\u00d734=\",\n  \"86\u00d789=\",\n  \"69\u00d731=\",\n  \"66\u00d728=\",\n  \"23\u00d743=\",\n  \"32\u00d716=\",\n  \"29\u00d717=\",\n  \"22\u00d714=\",\n  \"30\u00d771=\",\n  \"89\u00d767=\",\n  \"16\u00d724=\",\n  \"72\u00d726=\",\n  \"56\u00d723=\",\n  \"18\u00d713=\",\n  \"81\u00d796=\",\n  \"48\u00d765=\",\n  \"80\u00d763=\",\n  \"66\u00d726=\",\n  \"49\u00d772=\",\n  \"14\u00d799=\",\n  \"53\u00d721=\",\n  \"20\u00d712=\",\n  \"68\u00d790=\",\n  \"33\u00d760=\",\n  \"84\u00d785=\",\n  \"91\u00d720=\",\n  \"14\u00d766=\",\n  \"55\u00d750=\",\n  \"35\u00d773=\",\n  \"89\u00d728=\",\n  \"59\u00d775=\",\n  \"86\u00d779=\",\n  \"95\u00d732=\",\n  \"13\u00d711=\",\n  \"94\u00d793=\",\n  \"95\u00d751=\",\n  \"72\u00d713=\",\n  \"33\u00d728=\",\n  \"43\u00d769=\",\n  \"82\u00d797=\",\n  \"32\u00d711=\",\n  \"81\u00d725=\",\n  \"15\u00d778=\",\n  \"18\u00d745=\",\n  \"69\u00d753=\",\n  \"38\u00d736=\",\n  \"47\u00d733=\",\n  \"16\u00d756=\",\n  \"19\u00d767=\",\n  \"41\u00d716=\",\n  \"86\u00d740=\",\n  \"23\u00d767=\",\n  \"63\u00d720=\",\n  \"52\u00d712=\",\n  \"12\u00d797=\",\n  \"96\u00d747=\",\n  \"91\u00d773=\",\n  \"15\u00d781=\",\n  \"70\u00d771=\",\n  \"76\u00d770=\",\n  \"41\u00d771=\",\n  \"72\u00d735=\",\n  \"57\u00d785=\",\n  \"34\u00d713=\",\n  \"12\u00d745=\",\n  \"97\u00d752=\",\n)\n\n# --- Date paragraph (first paragraph in the body) ---------------------\n$dateOld = \"2023-06-30 Friday\"\n$dateNew = \"2023-07-01 Saturday\"\n$p1 = $d.Paragraphs(1)\n$p1Text = $p1.Range.Text.TrimEnd([char]13, [char]7)\nif ($p1Text -ne $dateOld) {\n    throw \"Paragraph 1: expected '$dateOld' but found '$p1Text'\"\n}\n$p1.Range.Text = $dateNew\n\n# --- Multiplication table (20 rows x 5 columns, row-major) -------------\n\n$t = $d.Tables(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\nif ($rows -ne 20 -or $cols -ne 5) {\n    throw \"Expected a 20x5 table, found ${rows}x${cols}\"\n}\nif ($OLD_VALUES.Count -ne 100 -or $NEW_VALUES.Count -ne 100) {\n    throw \"Expected 100 old/new cell values\"\n}\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cellText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        $expectedOld = $OLD_VALUES[$idx]\n        $expectedNew = $NEW_VALUES[$idx]\n        if ($cellText -ne $expectedOld) {\n            throw \"Cell ($r,$c) [index $idx]: expected '$expectedOld' but found '$cellText'\"\n        }\n        $cell.Range.Text = $expectedNew\n        $idx++\n    }\n}\n\nWrite-Output \"done\"\n"}
